$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("covid19_cases_switzerland")
$ws2 = $wb.Worksheets.Item("Quellen")

# --- covid19_cases_switzerland: add "yesterday"'s AG figure for 2020-03-16 (row 12) ---
$ws1.Range("B12").Value = 52

# --- Quellen: add four new cantons' sources ---

# AG
$ws2.Range("A7").Value = "AG"
$ws2.Range("B7").Value = "https://www.ag.ch/media/kanton_aargau/themen_1/coronavirus_1/200316_Coronavirus_Bundesrats_Entscheide-2.pdf"
[void]$ws2.Hyperlinks.Add($ws2.Range("B7"), "https://www.ag.ch/media/kanton_aargau/themen_1/coronavirus_1/200316_Coronavirus_Bundesrats_Entscheide-2.pdf")
$ws2.Range("B7").Style = "Hyperlink"
$ws2.Range("C7").Value = "@BachliMeyer"
$ws2.Range("C5").Copy() | Out-Null
[void]$ws2.Range("C7").PasteSpecial(-4122)

# GR
$ws2.Range("A8").Value = "GR"
$ws2.Range("B8").Value = "https://www.youtube.com/channel/UCEcqzK6vbCuIvxLiJCAMVuA"
[void]$ws2.Hyperlinks.Add($ws2.Range("B8"), "https://www.youtube.com/channel/UCEcqzK6vbCuIvxLiJCAMVuA")
$ws2.Range("B8").Style = "Hyperlink"

# UR
$ws2.Range("A10").Value = "UR"
$ws2.Range("B10").Value = "https://www.ur.ch/mmdirektionen/63458"
[void]$ws2.Hyperlinks.Add($ws2.Range("B10"), "https://www.ur.ch/mmdirektionen/63458")
$ws2.Range("B10").Style = "Hyperlink"
$ws2.Range("C10").Value = "@BachliMeyer"
$ws2.Range("C5").Copy() | Out-Null
[void]$ws2.Range("C10").PasteSpecial(-4122)

# AI
$ws2.Range("A9").Value = "AI"
$ws2.Range("B9").Value = "https://www.ai.ch/themen/gesundheit-alter-und-soziales/gesundheitsfoerderung-und-praevention/uebertragbare-krankheiten/coronavirus"
[void]$ws2.Hyperlinks.Add($ws2.Range("B9"), "https://www.ai.ch/themen/gesundheit-alter-und-soziales/gesundheitsfoerderung-und-praevention/uebertragbare-krankheiten/coronavirus")
$ws2.Range("B9").Style = "Hyperlink"
$ws2.Range("C9").Value = "@BachliMeyer"
$ws2.Range("C5").Copy() | Out-Null
[void]$ws2.Range("C9").PasteSpecial(-4122)

# --- view state: selection on sheet1, then switch active tab to Quellen with its own selection ---
[void]$ws1.Range("G23").Select()
[void]$ws2.Activate()
[void]$ws2.Range("C16").Select()
